$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update header row text (row 1)
$ws.Range("A1").Value = "Email address"
$ws.Range("B1").Value = "Password"
$ws.Range("C1").Value = "Result"

# Add a new (blank) row 5, formatted like the hyperlink cells above it (A2:A4)
$ws.Range("A5").Value = ""
$ws.Range("A2").Copy()
$ws.Range("A5").PasteSpecial(-4122)

# Update the selection shown when the workbook is opened
$ws.Range("A4").Select()
